$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.006.39'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.844.36'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.42'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.60'
$ws.Range('E8').Value = '  +6.12%  '
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0694'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0981'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.109.88'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.48'
$ws.Range('E13').Value = '  +4.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.842.45'
$ws.Range('E14').Value = '  +2.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.672'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.006.16'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.05'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('E19').Value = '  +0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.30'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.19'
$ws.Range('E21').Value = '  +2.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.77'
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.01'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.88'
$ws.Range('E26').Value = '  +2.04%  '
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.74'
$ws.Range('E29').Value = '  +11.95%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.93'
$ws.Range('E33').Value = '  -1.03%  '
$ws.Range('E34').Value = '  +23.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.97'
$ws.Range('E35').Value = '  +10.78%  '
$ws.Range('E36').Value = '  -3.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.754'
$ws.Range('E37').Value = '  +8.14%  '
$ws.Range('E38').Value = '  +9.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '89.92'
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('E40').Value = '  +3.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.344.09'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.58'
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.27'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.40'
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.76'
$ws.Range('E45').Value = '  +3.89%  '
$ws.Range('B46').Value = 'Gas'
$ws.Range('C46').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.92'
$ws.Range('E46').Value = '  +81.64%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0530'
$ws.Range('E47').Value = '  +3.56%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.30'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.023.26'
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.01'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0671'
$ws.Range('E51').Value = '  +0.08%  '
